{"js": "// Replace the 25 \"two-digit \u00f7 one-digit\" problems in the worksheet table\n// with a new set of problems, preserving the cell's run formatting\n// (font/size) untouched \u2014 only the <w:t> text content changes.\n//\n// The mapping below is positional/sequential: it lists the OLD problem\n// text followed by the NEW problem text, in document order. Every value\n// is unique in the document, so each pair unambiguously identifies one\n// paragraph.\nconst replacements = [\n  [\"68\u00f78=\", \"95\u00f73=\"],\n  [\"40\u00f77=\", \"77\u00f79=\"],\n  [\"68\u00f73=\", \"95\u00f78=\"],\n  [\"10\u00f75=\", \"94\u00f74=\"],\n  [\"71\u00f74=\", \"57\u00f78=\"],\n  [\"61\u00f72=\", \"79\u00f74=\"],\n  [\"69\u00f73=\", \"60\u00f79=\"],\n  [\"53\u00f76=\", \"25\u00f77=\"],\n  [\"41\u00f75=\", \"75\u00f78=\"],\n  [\"55\u00f73=\", \"28\u00f72=\"],\n  [\"85\u00f75=\", \"45\u00f73=\"],\n  [\"95\u00f73=\", \"22\u00f72=\"],\n  [\"69\u00f79=\", \"50\u00f73=\"],\n  [\"46\u00f72=\", \"17\u00f76=\"],\n  [\"93\u00f73=\", \"45\u00f75=\"],\n  [\"19\u00f79=\", \"30\u00f76=\"],\n  [\"97\u00f75=\", \"18\u00f75=\"],\n  [\"51\u00f75=\", \"58\u00f72=\"],\n  [\"79\u00f77=\", \"58\u00f79=\"],\n  [\"46\u00f77=\", \"84\u00f74=\"],\n  [\"33\u00f79=\", \"89\u00f74=\"],\n  [\"63\u00f75=\", \"55\u00f75=\"],\n  [\"16\u00f76=\", \"75\u00f78=\"],\n  [\"51\u00f76=\", \"60\u00f79=\"],\n  [\"53\u00f79=\", \"17\u00f78=\"],\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Walk the document's paragraphs in order and consume the replacement\n// list in lockstep: the Nth non-empty math-cell paragraph corresponds to\n// the Nth entry of `replacements` (this mirrors the diff, which edits the\n// table cells top-to-bottom, left-to-right, skipping the date line and\n// the blank filler cells).\nlet cursor = 0;\nfor (const paragraph of paragraphs.items) {\n  if (cursor >= replacements.length) break;\n  const [oldText, newText] = replacements[cursor];\n  if (paragraph.text === oldText) {\n    paragraph.insertText(newText, \"Replace\");\n    cursor++;\n  }\n}\n\nawait context.sync();\n\nif (cursor !== replacements.length) {\n  throw new Error(\n    `Only matched ${cursor} of ${replacements.length} expected cells`\n  );\n}\n", "ps1": "# Replace the 25 \"two-digit \u00f7 one-digit\" problems in the worksheet table\n# with a new set of problems, preserving each cell's run formatting\n# (font/size) untouched -- only the paragraph's text content changes.\n#\n# The list below is positional/sequential: each row is (OLD, NEW) in\n# document order. Every OLD value is unique in the document, so walking\n# the document's paragraphs top-to-bottom and consuming this list in\n# lockstep unambiguously identifies the right paragraph each time -- even\n# though some NEW values (e.g. \"95\u00f73=\") equal an OLD value used earlier\n# (e.g. \"68\u00f78=\" -> \"95\u00f73=\"), which would make a naive repeated\n# Find/Replace over the whole document ambiguous.\n$replacements = @(\n    @(\"68\u00f78=\", \"95\u00f73=\"),\n    @(\"40\u00f77=\", \"77\u00f79=\"),\n    @(\"68\u00f73=\", \"95\u00f78=\"),\n    @(\"10\u00f75=\", \"94\u00f74=\"),\n    @(\"71\u00f74=\", \"57\u00f78=\"),\n    @(\"61\u00f72=\", \"79\u00f74=\"),\n    @(\"69\u00f73=\", \"60\u00f79=\"),\n    @(\"53\u00f76=\", \"25\u00f77=\"),\n    @(\"41\u00f75=\", \"75\u00f78=\"),\n    @(\"55\u00f73=\", \"28\u00f72=\"),\n    @(\"85\u00f75=\", \"45\u00f73=\"),\n    @(\"95\u00f73=\", \"22\u00f72=\"),\n    @(\"69\u00f79=\", \"50\u00f73=\"),\n    @(\"46\u00f72=\", \"17\u00f76=\"),\n    @(\"93\u00f73=\", \"45\u00f75=\"),\n    @(\"19\u00f79=\", \"30\u00f76=\"),\n    @(\"97\u00f75=\", \"18\u00f75=\"),\n    @(\"51\u00f75=\", \"58\u00f72=\"),\n    @(\"79\u00f77=\", \"58\u00f79=\"),\n    @(\"46\u00f77=\", \"84\u00f74=\"),\n    @(\"33\u00f79=\", \"89\u00f74=\"),\n    @(\"63\u00f75=\", \"55\u00f75=\"),\n    @(\"16\u00f76=\", \"75\u00f78=\"),\n    @(\"51\u00f76=\", \"60\u00f79=\"),\n    @(\"53\u00f79=\", \"17\u00f78=\")\n)\n\n$d = $word.ActiveDocument\n\n$cursor = 0\n$total = $replacements.Length\n\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($cursor -ge $total) { break }\n\n    $para = $d.Paragraphs($i)\n    $rawText = $para.Range.Text\n    $text = $rawText.TrimEnd([char]13, [char]7)\n\n    $oldText = $replacements[$cursor][0]\n    $newText = $replacements[$cursor][1]\n\n    if ($text -eq $oldText) {\n        $para.Range.Text = $newText\n        $cursor++\n    }\n}\n\nif ($cursor -ne $total) {\n    throw \"Only matched $cursor of $total expected cells\"\n}\n\nWrite-Output \"done\"\n"}
